$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($ws, $addr, $text) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-CellText $ws "D2" '29.384.58'
Set-CellText $ws "E2" '  +0.04%  '

Set-CellText $ws "D3" '1.847.82'
Set-CellText $ws "E3" '  +0.08%  '

Set-CellText $ws "D4" '0.9996'
Set-CellText $ws "E4" '  -0.02%  '

Set-CellText $ws "D5" '239.98'
Set-CellText $ws "E5" '  -0.54%  '

Set-CellText $ws "E6" '  -0.50%  '

Set-CellText $ws "E7" '  -0.03%  '

Set-CellText $ws "D8" '0.07620'
Set-CellText $ws "E8" '  +0.79%  '

Set-CellText $ws "D9" '0.2928'
Set-CellText $ws "E9" '  -1.02%  '

Set-CellText $ws "D10" '24.49'
Set-CellText $ws "E10" '  -1.16%  '

Set-CellText $ws "D11" '0.07740'
Set-CellText $ws "E11" '  +0.00%  '

Set-CellText $ws "D12" '1.853.62'
Set-CellText $ws "E12" '  -6.62%  '

Set-CellText $ws "D13" '0.00001115'
Set-CellText $ws "E13" '  +11.98%  '

Set-CellText $ws "D14" '5.001'
Set-CellText $ws "E14" '  +0.13%  '

Set-CellText $ws "D15" '0.6779'
Set-CellText $ws "E15" '  -0.61%  '

Set-CellText $ws "D16" '83.68'
Set-CellText $ws "E16" '  +0.82%  '

Set-CellText $ws "D17" '2.107.44'
Set-CellText $ws "E17" '  -6.93%  '

Set-CellText $ws "D18" '6.175'

Set-CellText $ws "D19" '29.403.64'
Set-CellText $ws "E19" '  +0.01%  '

Set-CellText $ws "D20" '228.68'
Set-CellText $ws "E20" '  -0.69%  '

Set-CellText $ws "D21" '12.44'
Set-CellText $ws "E21" '  +0.13%  '

Set-CellText $ws "D22" '1.000'
Set-CellText $ws "E22" '  +0.05%  '

Set-CellText $ws "D23" '7.486'
Set-CellText $ws "E23" '  -0.82%  '

Set-CellText $ws "E24" '  +0.02%  '

Set-CellText $ws "D25" '157.22'
Set-CellText $ws "E25" '  +0.52%  '

Set-CellText $ws "D26" '0.1395'
Set-CellText $ws "E26" '  -0.18%  '

Set-CellText $ws "D27" '8.342'
Set-CellText $ws "E27" '  -0.49%  '

Set-CellText $ws "D28" '17.62'
Set-CellText $ws "E28" '  -0.26%  '

Set-CellText $ws "E29" '  -0.30%  '

Set-CellText $ws "D30" '1.300'
Set-CellText $ws "E30" '  +3.74%  '

Set-CellText $ws "D31" '0.05593'
Set-CellText $ws "E31" '  -2.08%  '

Set-CellText $ws "D32" '4.112'
Set-CellText $ws "E32" '  -0.28%  '

Set-CellText $ws "D33" '4.030'
Set-CellText $ws "E33" '  +0.56%  '

Set-CellText $ws "D34" '1.847'
Set-CellText $ws "E34" '  +0.28%  '

Set-CellText $ws "D35" '1.156'
Set-CellText $ws "E35" '  +0.21%  '

Set-CellText $ws "D36" '0.7106'
Set-CellText $ws "E36" '  -0.79%  '

Set-CellText $ws "D37" '2.584'
Set-CellText $ws "E37" '  -0.33%  '

Set-CellText $ws "D38" '1.238.86'
Set-CellText $ws "E38" '  -0.41%  '

Set-CellText $ws "B39" 'VeChain'
Set-CellText $ws "C39" 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-CellText $ws "D39" '0.01804'
Set-CellText $ws "E39" '  -0.19%  '

Set-CellText $ws "B40" 'MXToken'
Set-CellText $ws "C40" 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-CellText $ws "D40" '2.776'
Set-CellText $ws "E40" '  -0.90%  '

Set-CellText $ws "D41" '6.400'
Set-CellText $ws "E41" '  +5.04%  '

Set-CellText $ws "D42" '0.9048'
Set-CellText $ws "E42" '  +0.42%  '

Set-CellText $ws "D43" '1.0000'
Set-CellText $ws "E43" '  -0.05%  '

Set-CellText $ws "D44" '101.95'
Set-CellText $ws "E44" '  +0.08%  '

Set-CellText $ws "D45" '65.92'
Set-CellText $ws "E45" '  -0.33%  '

Set-CellText $ws "D46" '7.141'
Set-CellText $ws "E46" '  +1.32%  '

Set-CellText $ws "B47" 'TheSandbox'
Set-CellText $ws "C47" 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-CellText $ws "D47" '0.4013'
Set-CellText $ws "E47" '  -0.09%  '

Set-CellText $ws "B48" 'EnergySwap'
Set-CellText $ws "C48" 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-CellText $ws "D48" '9.007'
Set-CellText $ws "E48" '  -1.11%  '

Set-CellText $ws "B49" 'BabyDogeCoin'
Set-CellText $ws "C49" 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-CellText $ws "D49" '0.00000000115'
Set-CellText $ws "E49" '  -2.26%  '

Set-CellText $ws "D50" '1.682'
Set-CellText $ws "E50" '  -1.13%  '

Set-CellText $ws "E51" '  -0.36%  '
